# Update crypto price/volume values per latest data refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.119.96"
$ws.Range("E2").Value = "  +2.04%  "

$ws.Range("D3").Value = "2.363.12"
$ws.Range("E3").Value = "  +1.97%  "

$ws.Range("E4").Value = "  -0.54%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "543.76"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.66%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "136.21"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.85%  "

$ws.Range("E7").Value = "  +0.30%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.562"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +5.06%  "

$ws.Range("E9").Value = "  +1.58%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "5.60"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.99%  "

$ws.Range("E11").Value = "  -0.81%  "

$ws.Range("E12").Value = "  +1.02%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "24.00"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.92%  "

$ws.Range("D14").Value = "2.779.99"
$ws.Range("E14").Value = "  +1.61%  "

$ws.Range("D15").Value = "58.085.17"
$ws.Range("E15").Value = "  +1.40%  "

$ws.Range("E16").Value = "  +1.89%  "

$ws.Range("D17").Value = "2.354.18"
$ws.Range("E17").Value = "  +0.54%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "10.78"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +3.66%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "333.51"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.08%  "

$ws.Range("E20").Value = "  +2.48%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.80"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.35%  "

$ws.Range("E22").Value = "  +0.22%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "62.89"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.22%  "

$ws.Range("E24").Value = "  +0.43%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "8.54"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.83%  "

$ws.Range("E26").Value = "  +0.53%  "

$ws.Range("E27").Value = "  +2.33%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "172.99"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.22%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.76"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.43%  "

$ws.Range("D30").Value = "0.0₃0742"
$ws.Range("E30").Value = "  +2.50%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.18"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.07%  "

$ws.Range("E32").Value = "  +11.41%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "18.57"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.46%  "

$ws.Range("E34").Value = "  +0.04%  "

$ws.Range("E35").Value = "  +6.74%  "

$ws.Range("E36").Value = "  +0.74%  "

$ws.Range("E37").Value = "  +0.70%  "

$ws.Range("E38").Value = "  +3.98%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "39.44"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.55%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "145.58"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.54%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "293.93"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.57%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.380"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.34%  "

$ws.Range("E43").Value = "  +1.37%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0950"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.22%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "19.29"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +3.00%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0505"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.90%  "

$ws.Range("E47").Value = "  +1.07%  "

$ws.Range("E48").Value = "  +2.83%  "

$ws.Range("E49").Value = "  +0.73%  "

$ws.Range("E50").Value = "  +0.07%  "

$ws.Range("E51").Value = "  +0.38%  "
